$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "time_taken" header in F1, matching the formatting of the
# other header cells (bold font, thin border, centered) by copying E1's format.
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

$timestamps = @(
    "2021-10-05 13:41:02.277681",
    "2021-10-05 13:41:02.277691",
    "2021-10-05 13:41:02.277694",
    "2021-10-05 13:41:02.277697",
    "2021-10-05 13:41:02.277699",
    "2021-10-05 13:41:02.277702",
    "2021-10-05 13:41:02.277704",
    "2021-10-05 13:41:02.277707",
    "2021-10-05 13:41:02.277709",
    "2021-10-05 13:41:02.277712",
    "2021-10-05 13:41:02.277714",
    "2021-10-05 13:41:02.277716",
    "2021-10-05 13:41:02.277719",
    "2021-10-05 13:41:02.277721",
    "2021-10-05 13:41:02.277724",
    "2021-10-05 13:41:02.277726",
    "2021-10-05 13:41:02.277728",
    "2021-10-05 13:41:02.277731",
    "2021-10-05 13:41:02.277733",
    "2021-10-05 13:41:02.277736",
    "2021-10-05 13:41:02.277738",
    "2021-10-05 13:41:02.277740",
    "2021-10-05 13:41:02.277743",
    "2021-10-05 13:41:02.277745",
    "2021-10-05 13:41:02.277747",
    "2021-10-05 13:41:02.277750",
    "2021-10-05 13:41:02.277752",
    "2021-10-05 13:41:02.277754"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}

$excel.ActiveWorkbook.Save()
